$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "The American Journal of Gastroenterology"
$ws.Range("G2").Value = "https://openalex.org/S66441642"
$ws.Range("H2").Value = "Lippincott Williams & Wilkins"
$ws.Range("I2").Value = "0002-9270"
$ws.Range("V2").Value = "'FALSE"
$ws.Range("V2").Style = "Normal"
